$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.42873533333334
$ws.Range("H2").Value = 121.286206
$ws.Range("I2").Value = 0.08313576592793961
$ws.Range("J2").Value = 0.08313576592793961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.21127
$ws.Range("N2").Value = 0.63381
$ws.Range("O2").Value = 0.02468048274231428
$ws.Range("P2").Value = 0.02468048274231428
$ws.Range("Q2").Value = 8.541378913873332
$ws.Range("R2").Value = 76.87241022486
$ws.Range("S2").Value = 0.002051830836253593
$ws.Range("T2").Value = 0.002051830836253593

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.42873533333334
$ws.Range("H3").Value = 121.286206
$ws.Range("I3").Value = 0.08313576592793961
$ws.Range("J3").Value = 0.08313576592793961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.467027333333334
$ws.Range("N3").Value = 10.401082
$ws.Range("O3").Value = 0.405016842275123
$ws.Range("P3").Value = 0.405016842275123
$ws.Range("Q3").Value = 140.1675304527658
$ws.Range("R3").Value = 1261.507774074892
$ws.Range("S3").Value = 0.03367138539625786
$ws.Range("T3").Value = 0.03367138539625786

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.42873533333334
$ws.Range("H4").Value = 121.286206
$ws.Range("I4").Value = 0.08313576592793961
$ws.Range("J4").Value = 0.08313576592793961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.881908
$ws.Range("N4").Value = 14.645724
$ws.Range("O4").Value = 0.5703026749825627
$ws.Range("P4").Value = 0.5703026749825627
$ws.Range("Q4").Value = 197.3693664536827
$ws.Range("R4").Value = 1776.324298083144
$ws.Range("S4").Value = 0.04741254969542816
$ws.Range("T4").Value = 0.04741254969542816

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 412.4720866666667
$ws.Range("H5").Value = 1237.41626
$ws.Range("I5").Value = 0.84818836320749
$ws.Range("J5").Value = 0.84818836320749
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.21127
$ws.Range("N5").Value = 0.63381
$ws.Range("O5").Value = 0.02468048274231428
$ws.Range("P5").Value = 0.02468048274231428
$ws.Range("Q5").Value = 87.14297775006666
$ws.Range("R5").Value = 784.2867997505999
$ws.Range("S5").Value = 0.02093369826037425
$ws.Range("T5").Value = 0.02093369826037425

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 412.4720866666667
$ws.Range("H6").Value = 1237.41626
$ws.Range("I6").Value = 0.84818836320749
$ws.Range("J6").Value = 0.84818836320749
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.467027333333334
$ws.Range("N6").Value = 10.401082
$ws.Range("O6").Value = 0.405016842275123
$ws.Range("P6").Value = 0.405016842275123
$ws.Range("Q6").Value = 1430.051998710369
$ws.Range("R6").Value = 12870.46798839332
$ws.Range("S6").Value = 0.3435305725208028
$ws.Range("T6").Value = 0.3435305725208028

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 412.4720866666667
$ws.Range("H7").Value = 1237.41626
$ws.Range("I7").Value = 0.84818836320749
$ws.Range("J7").Value = 0.84818836320749
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.881908
$ws.Range("N7").Value = 14.645724
$ws.Range("O7").Value = 0.5703026749825627
$ws.Range("P7").Value = 0.5703026749825627
$ws.Range("Q7").Value = 2013.650779674693
$ws.Range("R7").Value = 18122.85701707224
$ws.Range("S7").Value = 0.483724092426313
$ws.Range("T7").Value = 0.483724092426313

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.396921
$ws.Range("H8").Value = 100.190763
$ws.Range("I8").Value = 0.0686758708645703
$ws.Range("J8").Value = 0.0686758708645703
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.21127
$ws.Range("N8").Value = 0.63381
$ws.Range("O8").Value = 0.02468048274231428
$ws.Range("P8").Value = 0.02468048274231428
$ws.Range("Q8").Value = 7.055767499669999
$ws.Range("R8").Value = 63.50190749702999
$ws.Range("S8").Value = 0.001694953645686431
$ws.Range("T8").Value = 0.001694953645686432

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.396921
$ws.Range("H9").Value = 100.190763
$ws.Range("I9").Value = 0.0686758708645703
$ws.Range("J9").Value = 0.0686758708645703
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.467027333333334
$ws.Range("N9").Value = 10.401082
$ws.Range("O9").Value = 0.405016842275123
$ws.Range("P9").Value = 0.405016842275123
$ws.Range("Q9").Value = 115.788037956174
$ws.Range("R9").Value = 1042.092341605566
$ws.Range("S9").Value = 0.02781488435806239
$ws.Range("T9").Value = 0.02781488435806239

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.396921
$ws.Range("H10").Value = 100.190763
$ws.Range("I10").Value = 0.0686758708645703
$ws.Range("J10").Value = 0.0686758708645703
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.881908
$ws.Range("N10").Value = 14.645724
$ws.Range("O10").Value = 0.5703026749825627
$ws.Range("P10").Value = 0.5703026749825627
$ws.Range("Q10").Value = 163.040695805268
$ws.Range("R10").Value = 1467.366262247412
$ws.Range("S10").Value = 0.03916603286082149
$ws.Range("T10").Value = 0.03916603286082149
